{"js": "// Replace the 25 division-problem cells in the single table with their new values.\n// The table has 20 rows x 5 columns; only rows 0, 4, 8, 12, 16 carry text (the\n// other rows are blank spacer rows). We address each target cell by its\n// (row, column) position and overwrite only the run text, via the paragraph's\n// range, so existing run/paragraph formatting (font, size, justification) is\n// preserved.\n\nconst newValues = [\n  // row 0\n  [\"88\u00f77=\", \"57\u00f76=\"],\n  [\"15\u00f77=\", \"65\u00f74=\"],\n  [\"27\u00f75=\", \"18\u00f72=\"],\n  [\"25\u00f79=\", \"27\u00f75=\"],\n  [\"54\u00f76=\", \"58\u00f73=\"],\n  // row 4\n  [\"35\u00f78=\", \"10\u00f72=\"],\n  [\"40\u00f77=\", \"88\u00f79=\"],\n  [\"98\u00f75=\", \"75\u00f75=\"],\n  [\"59\u00f72=\", \"22\u00f76=\"],\n  [\"30\u00f74=\", \"59\u00f72=\"],\n  // row 8\n  [\"89\u00f76=\", \"16\u00f78=\"],\n  [\"79\u00f76=\", \"55\u00f79=\"],\n  [\"30\u00f74=\", \"30\u00f78=\"],\n  [\"49\u00f73=\", \"46\u00f74=\"],\n  [\"21\u00f76=\", \"96\u00f72=\"],\n  // row 12\n  [\"11\u00f78=\", \"21\u00f75=\"],\n  [\"12\u00f72=\", \"59\u00f72=\"],\n  [\"37\u00f79=\", \"96\u00f77=\"],\n  [\"42\u00f77=\", \"57\u00f79=\"],\n  [\"59\u00f75=\", \"26\u00f72=\"],\n  // row 16\n  [\"47\u00f79=\", \"43\u00f75=\"],\n  [\"24\u00f75=\", \"80\u00f75=\"],\n  [\"13\u00f77=\", \"32\u00f79=\"],\n  [\"58\u00f75=\", \"80\u00f72=\"],\n  [\"62\u00f74=\", \"68\u00f74=\"],\n];\n\nconst contentRows = [0, 4, 8, 12, 16];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Collect the target paragraph ranges (first paragraph of each target cell).\nconst paragraphsByCell = [];\nfor (const row of contentRows) {\n  for (let col = 0; col < 5; col++) {\n    const cell = table.getCell(row, col);\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    paragraphsByCell.push(paragraphs);\n  }\n}\nawait context.sync();\n\nconst ranges = paragraphsByCell.map((p) => p.items[0].getRange());\nranges.forEach((r) => r.load(\"text\"));\nawait context.sync();\n\n// Replace each paragraph's text in place (via its Range), which keeps the\n// run's existing formatting (font, size) and the paragraph's justification\n// intact \u2014 only the literal characters change, exactly like the diff.\nfor (let idx = 0; idx < ranges.length; idx++) {\n  const [before, after] = newValues[idx];\n  if (ranges[idx].text !== before) {\n    throw new Error(\n      `Unexpected text in cell ${idx}: expected \"${before}\", found \"${ranges[idx].text}\"`\n    );\n  }\n  ranges[idx].insertText(after, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem cells in the single table with their new\n# values. The table has 20 rows x 5 columns; only rows 1, 5, 9, 13, 17\n# (1-indexed, as Word COM addresses them) carry text \u2014 the other rows are\n# blank spacer rows. Each target cell's Range.Text is overwritten directly\n# (not the whole cell deleted/recreated), which preserves the existing\n# run/paragraph formatting (font, size, justification) already present.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Each inner array is one content row of the table: five \"before=>after\"\n# pairs, left to right.\n$rows = @(\n    @(\n        @{ Before = \"88\u00f77=\"; After = \"57\u00f76=\" },\n        @{ Before = \"15\u00f77=\"; After = \"65\u00f74=\" },\n        @{ Before = \"27\u00f75=\"; After = \"18\u00f72=\" },\n        @{ Before = \"25\u00f79=\"; After = \"27\u00f75=\" },\n        @{ Before = \"54\u00f76=\"; After = \"58\u00f73=\" }\n    ),\n    @(\n        @{ Before = \"35\u00f78=\"; After = \"10\u00f72=\" },\n        @{ Before = \"40\u00f77=\"; After = \"88\u00f79=\" },\n        @{ Before = \"98\u00f75=\"; After = \"75\u00f75=\" },\n        @{ Before = \"59\u00f72=\"; After = \"22\u00f76=\" },\n        @{ Before = \"30\u00f74=\"; After = \"59\u00f72=\" }\n    ),\n    @(\n        @{ Before = \"89\u00f76=\"; After = \"16\u00f78=\" },\n        @{ Before = \"79\u00f76=\"; After = \"55\u00f79=\" },\n        @{ Before = \"30\u00f74=\"; After = \"30\u00f78=\" },\n        @{ Before = \"49\u00f73=\"; After = \"46\u00f74=\" },\n        @{ Before = \"21\u00f76=\"; After = \"96\u00f72=\" }\n    ),\n    @(\n        @{ Before = \"11\u00f78=\"; After = \"21\u00f75=\" },\n        @{ Before = \"12\u00f72=\"; After = \"59\u00f72=\" },\n        @{ Before = \"37\u00f79=\"; After = \"96\u00f77=\" },\n        @{ Before = \"42\u00f77=\"; After = \"57\u00f79=\" },\n        @{ Before = \"59\u00f75=\"; After = \"26\u00f72=\" }\n    ),\n    @(\n        @{ Before = \"47\u00f79=\"; After = \"43\u00f75=\" },\n        @{ Before = \"24\u00f75=\"; After = \"80\u00f75=\" },\n        @{ Before = \"13\u00f77=\"; After = \"32\u00f79=\" },\n        @{ Before = \"58\u00f75=\"; After = \"80\u00f72=\" },\n        @{ Before = \"62\u00f74=\"; After = \"68\u00f74=\" }\n    )\n)\n\n$contentRowIndexes = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $contentRowIndexes.Length; $i++) {\n    $tableRow = $contentRowIndexes[$i]\n    $cellDefs = $rows[$i]\n    for ($col = 1; $col -le 5; $col++) {\n        $def = $cellDefs[$col - 1]\n        $cell = $t.Cell($tableRow, $col)\n        $cellRange = $cell.Range\n        # Cell.Range.Text includes the trailing end-of-cell marker\n        # (CR + cell mark); strip it before comparing to the expected value.\n        $currentText = $cellRange.Text.TrimEnd([char]13, [char]7)\n        if ($currentText -ne $def.Before) {\n            throw \"Unexpected text in row $tableRow col $col`: expected '$($def.Before)', found '$currentText'\"\n        }\n        $cellRange.Text = $def.After\n    }\n}\n"}
